$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.387.67"
$ws.Range("E2").Value = "'  +4.43%  "
$ws.Range("D3").Value = "'3.634.28"
$ws.Range("E3").Value = "'  +4.16%  "
$ws.Range("D5").Value = "'593.55"
$ws.Range("E5").Value = "'  +1.20%  "
$ws.Range("D6").Value = "'195.34"
$ws.Range("E6").Value = "'  +4.92%  "
$ws.Range("D7").Value = "'0.645"
$ws.Range("E7").Value = "'  +1.57%  "
$ws.Range("D8").Value = "'3.628.10"
$ws.Range("E8").Value = "'  +4.16%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "'  +0.04%  "
$ws.Range("E10").Value = "'  +3.60%  "
$ws.Range("E11").Value = "'  +2.63%  "
$ws.Range("D12").Value = "'58.81"
$ws.Range("E12").Value = "'  +4.24%  "
$ws.Range("D13").Value = "'0.0000292"
$ws.Range("E13").Value = "'  +4.01%  "
$ws.Range("D14").Value = "'9.95"
$ws.Range("E14").Value = "'  +5.00%  "
$ws.Range("D15").Value = "'4.212.15"
$ws.Range("E15").Value = "'  +4.09%  "
$ws.Range("D16").Value = "'19.89"
$ws.Range("E16").Value = "'  +4.91%  "
$ws.Range("D17").Value = "'3.624.41"
$ws.Range("E17").Value = "'  +3.69%  "
$ws.Range("D18").Value = "'70.340.29"
$ws.Range("E18").Value = "'  +4.37%  "
$ws.Range("D19").Value = "'12.71"
$ws.Range("E19").Value = "'  +4.05%  "
$ws.Range("E20").Value = "'  +1.46%  "
$ws.Range("E21").Value = "'  +4.13%  "
$ws.Range("D22").Value = "'487.62"
$ws.Range("E22").Value = "'  -0.53%  "
$ws.Range("D23").Value = "'19.29"
$ws.Range("E23").Value = "'  +15.15%  "
$ws.Range("E24").Value = "'  -1.67%  "
$ws.Range("E25").Value = "'  +0.75%  "
$ws.Range("D26").Value = "'91.37"
$ws.Range("E26").Value = "'  +1.32%  "
$ws.Range("D27").Value = "'3.17"
$ws.Range("E27").Value = "'  +7.21%  "
$ws.Range("D28").Value = "'11.49"
$ws.Range("E28").Value = "'  +4.35%  "
$ws.Range("D29").Value = "'9.62"
$ws.Range("E29").Value = "'  +4.88%  "
$ws.Range("B30").Value = "'EthereumClassic"
$ws.Range("C30").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'32.98"
$ws.Range("E30").Value = "'  +4.48%  "
$ws.Range("B31").Value = "'NEARProtocol"
$ws.Range("C31").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.92"
$ws.Range("E31").Value = "'  +9.75%  "
$ws.Range("E32").Value = "'  +7.86%  "
$ws.Range("D33").Value = "'627.65"
$ws.Range("E33").Value = "'  +5.18%  "
$ws.Range("E34").Value = "'  +4.21%  "
$ws.Range("D35").Value = "'66.07"
$ws.Range("E35").Value = "'  +2.95%  "
$ws.Range("D36").Value = "'41.22"
$ws.Range("E36").Value = "'  +12.21%  "
$ws.Range("D37").Value = "'0.413"
$ws.Range("E37").Value = "'  +6.39%  "
$ws.Range("D38").Value = "'0.0₃0823"
$ws.Range("E38").Value = "'  +7.13%  "
$ws.Range("E39").Value = "'  +0.03%  "
$ws.Range("D40").Value = "'0.147"
$ws.Range("E40").Value = "'  -2.32%  "
$ws.Range("E41").Value = "'  +0.79%  "
$ws.Range("D42").Value = "'3.299.10"
$ws.Range("E42").Value = "'  +0.97%  "
$ws.Range("D43").Value = "'3.16"
$ws.Range("E43").Value = "'  +7.74%  "
$ws.Range("D44").Value = "'2.83"
$ws.Range("E44").Value = "'  +11.05%  "
$ws.Range("D45").Value = "'0.0455"
$ws.Range("E45").Value = "'  +5.38%  "
$ws.Range("D46").Value = "'2.87"
$ws.Range("E46").Value = "'  +2.79%  "
$ws.Range("D47").Value = "'3.31"
$ws.Range("E47").Value = "'  +1.91%  "
$ws.Range("E48").Value = "'  +2.17%  "
$ws.Range("D49").Value = "'9.19"
$ws.Range("E49").Value = "'  +4.50%  "
$ws.Range("D50").Value = "'3.34"
$ws.Range("E50").Value = "'  +2.00%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "'  -0.06%  "
